# Auto-generated: update cryptocurrency price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.897.94"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "1.676.90"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'219.55"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'29.10"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "'0.265"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").Value = "'0.0640"
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.920.52"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "1.669.72"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").Value = "'0.603"
$ws.Range("E14").Value = "  +6.71%  "
$ws.Range("D15").Value = "'10.03"
$ws.Range("E15").Value = "  +7.14%  "
$ws.Range("E16").Value = "  +6.80%  "
$ws.Range("D17").Value = "30.899.84"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "'66.03"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'246.76"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").Value = "'9.95"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'159.09"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "'15.81"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").Value = "'6.67"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").Value = "'3.49"
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.32"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.518.44"
$ws.Range("E34").Value = "  +6.39%  "
$ws.Range("D35").Value = "'1.74"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("D36").Value = "'84.39"
$ws.Range("E36").Value = "  +12.50%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +8.76%  "
$ws.Range("D39").Value = "'0.0179"
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'2.65"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").Value = "'0.837"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "'0.0503"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("D48").Value = "'51.20"
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("D49").Value = "1.812.53"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").Value = "0.0₆0119"
$ws.Range("E50").Value = "  +7.43%  "
$ws.Range("D51").Value = "'93.07"
$ws.Range("E51").Value = "  +1.65%  "
